$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking text cells in column D as Text so they are not
# coerced into numbers by Excels automatic type detection.
$textCells = @("D15","D16","D17","D19","D20","D21")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("A2").Value = "아이빔테크놀로지"
$ws.Range("B2").Value = "2024.07.15~07.19"
$ws.Range("C2").Value = "7,300~8,500"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 16308
$ws.Range("F2").Value = "삼성증권"

$ws.Range("A3").Value = "피앤에스미캐닉스"
$ws.Range("B3").Value = "2024.07.11~07.17"
$ws.Range("C3").Value = "14,000~17,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 18900
$ws.Range("F3").Value = "키움증권"

$ws.Range("A4").Value = "티디에스팜"
$ws.Range("B4").Value = "2024.07.10~07.16"
$ws.Range("C4").Value = "9,500~10,700"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 9500
$ws.Range("F4").Value = "한국투자증권"

$ws.Range("A5").Value = "케이쓰리아이"
$ws.Range("B5").Value = "2024.07.10~07.16"
$ws.Range("C5").Value = "12,500~15,500"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 17500
$ws.Range("F5").Value = "하나증권"

$ws.Range("A6").Value = "NH스팩31호"
$ws.Range("B6").Value = "2024.07.09~07.10"
$ws.Range("C6").Value = "2,000~2,000"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 12000
$ws.Range("F6").Value = "NH투자증권"

$ws.Range("A7").Value = "SK증권스팩13호"
$ws.Range("B7").Value = "2024.07.09~07.10"
$ws.Range("C7").Value = "2,000~2,000"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 8000
$ws.Range("F7").Value = "SK증권"

$ws.Range("A8").Value = "산일전기(유가)"
$ws.Range("B8").Value = "2024.07.09~07.15"
$ws.Range("C8").Value = "24,000~30,000"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = 182400
$ws.Range("F8").Value = "미래에셋증권,삼성증권"

$ws.Range("A9").Value = "뱅크웨어글로벌"
$ws.Range("B9").Value = "2024.07.08~07.12"
$ws.Range("C9").Value = "16,000~19,000"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = 22400
$ws.Range("F9").Value = "미래에셋증권"

$ws.Range("A10").Value = "키움스팩9호"
$ws.Range("B10").Value = "2024.07.02~07.03"
$ws.Range("C10").Value = "2,000~2,000"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = 6000
$ws.Range("F10").Value = "키움증권"

$ws.Range("A11").Value = "넥스트바이오메디컬"
$ws.Range("B11").Value = "2024.07.01~07.05"
$ws.Range("C11").Value = "24,000~29,000"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = 24000
$ws.Range("F11").Value = "한국투자증권"

$ws.Range("A12").Value = "이베스트스팩6호"
$ws.Range("B12").Value = "2024.06.27~06.28"
$ws.Range("C12").Value = "2,000~2,000"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = 8000
$ws.Range("F12").Value = "엘에스증권"

$ws.Range("A13").Value = "엑셀세라퓨틱스"
$ws.Range("B13").Value = "2024.06.24~06.28"
$ws.Range("C13").Value = "6,200~7,700"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = 10032
$ws.Range("F13").Value = "대신증권"

$ws.Range("A14").Value = "이엔셀"
$ws.Range("B14").Value = "2024.06.17~06.21"
$ws.Range("C14").Value = "13,600~15,300"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = 21308
$ws.Range("F14").Value = "NH투자증권"

$ws.Range("A15").Value = "하스"
$ws.Range("B15").Value = "2024.06.13~06.19"
$ws.Range("C15").Value = "9,000~12,000"
$ws.Range("D15").Value = "16000"
$ws.Range("E15").Value = 16290
$ws.Range("F15").Value = "삼성증권"

$ws.Range("A16").Value = "에이치브이엠(구.한국진공야금)"
$ws.Range("B16").Value = "2024.06.11~06.17"
$ws.Range("C16").Value = "11,000~14,200"
$ws.Range("D16").Value = "18000"
$ws.Range("E16").Value = 26400
$ws.Range("F16").Value = "NH투자증권"

$ws.Range("A17").Value = "이노스페이스"
$ws.Range("B17").Value = "2024.06.11~06.17"
$ws.Range("C17").Value = "36,400~43,300"
$ws.Range("D17").Value = "43300"
$ws.Range("E17").Value = 48412
$ws.Range("F17").Value = "미래에셋증권,신한투자증권"

$ws.Range("A18").Value = "한국스팩15호"
$ws.Range("B18").Value = "2024.06.10~06.11"
$ws.Range("C18").Value = "2,000~2,000"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = 12500
$ws.Range("F18").Value = "한국투자증권"

$ws.Range("A19").Value = "하이젠알앤엠"
$ws.Range("B19").Value = "2024.06.07~06.13"
$ws.Range("C19").Value = "4,500~5,500"
$ws.Range("D19").Value = "7000"
$ws.Range("E19").Value = 15300
$ws.Range("F19").Value = "한국투자증권"

$ws.Range("A20").Value = "미래에셋비전스팩6호"
$ws.Range("B20").Value = "2024.06.05~06.07"
$ws.Range("C20").Value = "2,000~2,000"
$ws.Range("D20").Value = "2000"
$ws.Range("E20").Value = 12900
$ws.Range("F20").Value = "미래에셋증권"

$ws.Range("A21").Value = "KB스팩29호"
$ws.Range("B21").Value = "2024.06.04~06.05"
$ws.Range("C21").Value = "2,000~2,000"
$ws.Range("D21").Value = "2000"
$ws.Range("E21").Value = 12000
$ws.Range("F21").Value = "KB증권"
